$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.60769999999999
$ws.Range("C4").Value = -11.79960000000001
$ws.Range("A6").Value = -22.66820000000002
$ws.Range("D6").Value = -8.153599999999997
$ws.Range("A7").Value = -18.86999999999999
$ws.Range("D7").Value = -8.282199999999992
$ws.Range("A8").Value = -21.9364
$ws.Range("C8").Value = -13.212
$ws.Range("D8").Value = -8.866799999999996
$ws.Range("C9").Value = -10.24470000000001
$ws.Range("D10").Value = -7.780999999999997
$ws.Range("C12").Value = -10.25669999999999
$ws.Range("D13").Value = -8.531899999999991
$ws.Range("D14").Value = -7.2779
$ws.Range("A16").Value = -22.0177
$ws.Range("D16").Value = -8.544699999999995
$ws.Range("C17").Value = -14.33439999999998
$ws.Range("C18").Value = -12.55549999999999
$ws.Range("C19").Value = -11.2757
$ws.Range("A20").Value = -20.09429999999998
$ws.Range("C20").Value = -12.4408
$ws.Range("A21").Value = -20.44359999999998
$ws.Range("C26").Value = -12.03030000000002
$ws.Range("A28").Value = -22.0277
$ws.Range("A29").Value = -21.34939999999997
$ws.Range("A30").Value = -21.49320000000002
$ws.Range("D30").Value = -6.946099999999992
$ws.Range("C31").Value = -13.00090000000001
$ws.Range("A32").Value = -21.35880000000002
$ws.Range("D37").Value = -8.676499999999997
$ws.Range("C39").Value = -11.86749999999999
$ws.Range("A40").Value = -20.32620000000001
$ws.Range("C40").Value = -12.19310000000001
$ws.Range("D40").Value = -7.369699999999999
$ws.Range("C41").Value = -12.14629999999999
$ws.Range("C42").Value = -11.7803
$ws.Range("C43").Value = -12.26909999999999
$ws.Range("D44").Value = -6.593700000000007
$ws.Range("A46").Value = -21.82730000000001
$ws.Range("C47").Value = -12.20019999999999
$ws.Range("C48").Value = -11.458
$ws.Range("A51").Value = -21.69289999999999
$ws.Range("A52").Value = -22.2427
$ws.Range("C54").Value = -13.3123
$ws.Range("A57").Value = -22.44840000000002
$ws.Range("A59").Value = -22.24470000000002
$ws.Range("A62").Value = -21.9822
$ws.Range("C62").Value = -13.2281
$ws.Range("C63").Value = -10.3647
$ws.Range("C64").Value = -10.21929999999999
$ws.Range("A66").Value = -21.4373
$ws.Range("D70").Value = -6.652299999999998
$ws.Range("A73").Value = -20.44399999999997
$ws.Range("A74").Value = -21.49629999999997
$ws.Range("C76").Value = -12.17100000000001
$ws.Range("A77").Value = -20.541
$ws.Range("C81").Value = -13.06419999999999
$ws.Range("C84").Value = -13.2324
$ws.Range("C89").Value = -14.77039999999999
$ws.Range("D89").Value = -8.457799999999999
$ws.Range("D91").Value = -8.304499999999996
$ws.Range("A92").Value = -21.6526
$ws.Range("D93").Value = -6.381499999999994
$ws.Range("C94").Value = -10.53459999999999
$ws.Range("D98").Value = -7.060800000000003
$ws.Range("A100").Value = -22.07220000000002
